$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$adminJson = "{`n  ""username"" : ""admin01"",`n  ""userPassword"" : ""12345678""`n}"
$studentJson = "{`n  ""username"" : ""student"",`n  ""userPassword"" : ""12345678""`n}"

# --- Clear the old "pass" column (B) on both login sheets ---
$ws1.Range("B1:B2").ClearContents()
$ws2.Range("B1:B2").ClearContents()

# --- New header + JSON body values ---
$ws1.Range("A1").Value = "data"
$ws1.Range("A2").Value = $adminJson

$ws2.Range("A1").Value = "data"
$ws2.Range("A2").Value = $studentJson

# --- Formatting: data (row2) cells get wrap text + vertical-center, on a
#     distinct font object; build it once on sheet1 then copy the format
#     across so both sheets share the very same style index. ---
$dataCell = $ws1.Range("A2")
$dataCell.Font.ThemeColor = 1
$dataCell.WrapText = $true
$dataCell.VerticalAlignment = -4108
$dataCell.Copy()
$ws2.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Formatting: header (row1) cells get horizontal-center ---
$ws1.Range("A1").HorizontalAlignment = -4108

$hdr2 = $ws2.Range("A1")
$hdr2.Font.ThemeColor = 1
$hdr2.HorizontalAlignment = -4108

# --- Row heights for the JSON body rows ---
$ws1.Rows.Item(2).RowHeight = 60
$ws2.Rows.Item(2).RowHeight = 60

# --- Column widths ---
$ws1.Columns.Item(1).ColumnWidth = 28.65
$ws2.Columns.Item(1).ColumnWidth = 27.85

# --- Selections ---
[void]$ws1.Range("A2").Select()
[void]$ws2.Range("C8").Select()

# --- Page setup (sheet1 only, as in the target) ---
$ws1.PageSetup.Orientation = 1
$ws1.PageSetup.PaperSize = 9
